$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blog")

# Row 10 (Documentation): add a description of the change plus the date it was made
$ws.Range("D10").Value = "Added simple documentation (must expand later)https://github.com/Hanandrof/Group-3-Repository/blob/main/Documentation/Blog_Documentation.pdf"
$ws.Range("D10").WrapText = $true
$ws.Range("E10").Value = "4/27/2021"
$ws.Range("E10").NumberFormat = "mm-dd-yy"

# Row 10 grew taller to accommodate the wrapped description text
$ws.Rows.Item(10).RowHeight = 90

# Row 9 (Not needlessly complex): record that this isn't working
$ws.Range("D9").Value = "Not working"

# Row 7 (Wordpress Plugins Interactivity): update the issue description and the
# recommendation to reference "Paid Memberships Pro" instead of "Buddy Press Groups"
$ws.Range("C7").Value = "Look into the documentation of Buddy Press and Paid Memberships Pro and make sure they work together"
$ws.Range("B7").Value = "Interaction between Paid Memberships Pro and Buddy Press is not working. Groups is trying to give role permissions while Buddy Press is not accepting them."

# Row 7, column D: add a note on the Github repo / description of change, wrapped
$ws.Range("D7").Value = "I believe it has something to do with the multiple user plugins we had installed and it messed with the overall website"
$ws.Range("D7").WrapText = $true

# Update the active selection to match where the editor last left off
$ws.Range("E7").Select() | Out-Null
